# This script applies the "Saldo" workbook update:
#   1) Adds a new account (005681354 / MATHEUS / 25400) in the row just
#      above account 005064129 (THIAGO), keeping the sheet's descending
#      sort order by Saldo.
#   2) Updates account 004452597 (LARA): its balance moved from 52.76 to
#      7775.61, which also moves its row from near the bottom of the
#      sheet (sorted by Saldo) up to just above account 004207278
#      (CESAR). We implement this as: delete the old LARA row (the one
#      whose current balance is 52.76), then insert a new LARA row in
#      its new sorted position (just above CESAR) with the new balance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$missing = [System.Reflection.Missing]::Value

# --- Step 1: insert MATHEUS row right above THIAGO (005064129) ---
$thiagoRow = $ws.Range("A1").EntireColumn.Find("005064129", $missing, $missing, 1).Row
$ws.Rows.Item($thiagoRow).Insert()

$ws.Range("A" + $thiagoRow).NumberFormat = "@"
$ws.Range("A" + $thiagoRow).Value = "005681354"
$ws.Range("B" + $thiagoRow).Value = "MATHEUS"
$ws.Range("C" + $thiagoRow).Value = 25400

# --- Step 2: delete the old LARA row (account 004452597, balance 52.76) ---
$oldLaraRow = $ws.Range("A1").EntireColumn.Find("004452597", $missing, $missing, 1).Row
$ws.Rows.Item($oldLaraRow).Delete()

# --- Step 3: insert the updated LARA row right above CESAR (004207278) ---
$cesarRow = $ws.Range("A1").EntireColumn.Find("004207278", $missing, $missing, 1).Row
$ws.Rows.Item($cesarRow).Insert()

$ws.Range("A" + $cesarRow).NumberFormat = "@"
$ws.Range("A" + $cesarRow).Value = "004452597"
$ws.Range("B" + $cesarRow).Value = "LARA"
$ws.Range("C" + $cesarRow).Value = 7775.61

Write-Output "MATHEUS inserted at row $thiagoRow"
Write-Output "Old LARA deleted at row $oldLaraRow"
Write-Output "New LARA inserted at row $cesarRow"
